# Auto-generated script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "30.389.31"
Set-TextCell $ws.Range("E2") "  +0.19%  "

# Row 3
Set-TextCell $ws.Range("D3") "1.939.15"
Set-TextCell $ws.Range("E3") "  +0.24%  "

# Row 4
Set-TextCell $ws.Range("D4") "0.9998"
Set-TextCell $ws.Range("E4") "  -0.11%  "

# Row 5
Set-TextCell $ws.Range("D5") "0.7702"
Set-TextCell $ws.Range("E5") "  +8.58%  "

# Row 6
Set-TextCell $ws.Range("D6") "247.85"
Set-TextCell $ws.Range("E6") "  -1.11%  "

# Row 7
Set-TextCell $ws.Range("D7") "0.9997"
Set-TextCell $ws.Range("E7") "  -0.08%  "

# Row 8
Set-TextCell $ws.Range("E8") "  +1.28%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.3227"
Set-TextCell $ws.Range("E9") "  -2.36%  "

# Row 10
Set-TextCell $ws.Range("D10") "0.07122"
Set-TextCell $ws.Range("E10") "  -2.42%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.7872"
Set-TextCell $ws.Range("E11") "  -2.27%  "

# Row 12
Set-TextCell $ws.Range("D12") "0.08030"
Set-TextCell $ws.Range("E12") "  -0.48%  "

# Row 13
Set-TextCell $ws.Range("D13") "1.938.10"
Set-TextCell $ws.Range("E13") "  +0.14%  "

# Row 14
Set-TextCell $ws.Range("D14") "5.389"
Set-TextCell $ws.Range("E14") "  -1.74%  "

# Row 15
Set-TextCell $ws.Range("D15") "95.16"
Set-TextCell $ws.Range("E15") "  +0.55%  "

# Row 16
Set-TextCell $ws.Range("D16") "14.59"
Set-TextCell $ws.Range("E16") "  -3.56%  "

# Row 17
Set-TextCell $ws.Range("D17") "30.393.64"
Set-TextCell $ws.Range("E17") "  +0.20%  "

# Row 18
Set-TextCell $ws.Range("D18") "256.07"
Set-TextCell $ws.Range("E18") "  +1.05%  "

# Row 19
Set-TextCell $ws.Range("D19") "0.000008022"
Set-TextCell $ws.Range("E19") "  -2.30%  "

# Row 20
Set-TextCell $ws.Range("D20") "5.849"
Set-TextCell $ws.Range("E20") "  +0.94%  "

# Row 21
Set-TextCell $ws.Range("D21") "2.195.31"
Set-TextCell $ws.Range("E21") "  +0.23%  "

# Row 22
Set-TextCell $ws.Range("D22") "0.9998"
Set-TextCell $ws.Range("E22") "  -0.08%  "

# Row 23
Set-TextCell $ws.Range("D23") "0.9994"
Set-TextCell $ws.Range("E23") "  -0.15%  "

# Row 24
Set-TextCell $ws.Range("D24") "6.768"
Set-TextCell $ws.Range("E24") "  -3.07%  "

# Row 25
Set-TextCell $ws.Range("D25") "9.627"
Set-TextCell $ws.Range("E25") "  -1.31%  "

# Row 26
Set-TextCell $ws.Range("D26") "163.91"
Set-TextCell $ws.Range("E26") "  -0.65%  "

# Row 27
Set-TextCell $ws.Range("B27") "Stellar"
Set-TextCell $ws.Range("C27") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D27") "0.1347"
Set-TextCell $ws.Range("E27") "  +4.93%  "

# Row 28
Set-TextCell $ws.Range("B28") "EthereumClassic"
Set-TextCell $ws.Range("C28") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws.Range("D28") "19.16"
Set-TextCell $ws.Range("E28") "  -1.13%  "

# Row 29
Set-TextCell $ws.Range("D29") "2.306"
Set-TextCell $ws.Range("E29") "  -1.71%  "

# Row 30
Set-TextCell $ws.Range("E30") "  +1.22%  "

# Row 31
Set-TextCell $ws.Range("D31") "1.527"
Set-TextCell $ws.Range("E31") "  -0.87%  "

# Row 32
Set-TextCell $ws.Range("D32") "4.444"
Set-TextCell $ws.Range("E32") "  +0.36%  "

# Row 33
Set-TextCell $ws.Range("D33") "4.156"
Set-TextCell $ws.Range("E33") "  -0.47%  "

# Row 34
Set-TextCell $ws.Range("D34") "0.05204"
Set-TextCell $ws.Range("E34") "  -0.03%  "

# Row 35
Set-TextCell $ws.Range("D35") "1.286"
Set-TextCell $ws.Range("E35") "  +1.69%  "

# Row 36
Set-TextCell $ws.Range("D36") "0.7525"
Set-TextCell $ws.Range("E36") "  +0.56%  "

# Row 37
Set-TextCell $ws.Range("D37") "2.774"
Set-TextCell $ws.Range("E37") "  -0.66%  "

# Row 38
Set-TextCell $ws.Range("D38") "0.01976"
Set-TextCell $ws.Range("E38") "  +0.35%  "

# Row 39
Set-TextCell $ws.Range("D39") "2.813"
Set-TextCell $ws.Range("E39") "  +0.13%  "

# Row 40
Set-TextCell $ws.Range("D40") "79.14"
Set-TextCell $ws.Range("E40") "  +0.21%  "

# Row 41
Set-TextCell $ws.Range("D41") "6.534"
Set-TextCell $ws.Range("E41") "  +1.80%  "

# Row 42
Set-TextCell $ws.Range("D42") "0.4538"
Set-TextCell $ws.Range("E42") "  +0.09%  "

# Row 43
Set-TextCell $ws.Range("D43") "1.987"
Set-TextCell $ws.Range("E43") "  -1.65%  "

# Row 44
Set-TextCell $ws.Range("E44") "  +0.01%  "

# Row 45
Set-TextCell $ws.Range("D45") "0.8374"
Set-TextCell $ws.Range("E45") "  -0.90%  "

# Row 46
Set-TextCell $ws.Range("D46") "101.56"
Set-TextCell $ws.Range("E46") "  -0.12%  "

# Row 47
Set-TextCell $ws.Range("D47") "9.861"
Set-TextCell $ws.Range("E47") "  +1.15%  "

# Row 48
Set-TextCell $ws.Range("D48") "7.521"
Set-TextCell $ws.Range("E48") "  +0.87%  "

# Row 49
Set-TextCell $ws.Range("B49") "Elrond"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell $ws.Range("D49") "37.50"
Set-TextCell $ws.Range("E49") "  +1.96%  "

# Row 50
Set-TextCell $ws.Range("B50") "Maker"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Range("D50") "985.46"
Set-TextCell $ws.Range("E50") "  +11.48%  "

# Row 51
Set-TextCell $ws.Range("E51") "  -0.12%  "
